# "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
#
# The "Periodo Mora" / "Valor Mora" table (rows 16-27, columns E & F) gets
# its period list reversed (most-recent period 2303 now on top, oldest 2204
# at the bottom) - old account statements are dropped off the top and new
# ones appended, so the whole list order flips. The "Valor Mora" column
# keeps riding along with its row's period (only the two rows that actually
# had a different value - 37333 vs 40000 - show up changed once reversed).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periods = @("2303","2302","2301","2212","2211","2210","2209","2208","2207","2206","2205","2204")
$valores = @(37333,40000,40000,40000,40000,40000,40000,40000,40000,40000,40000,40000)

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periods[$i]
    $ws.Range("F$row").Value = $valores[$i]
}

# Column widths were re-measured by Excel (bestFit) after the edit; nudge
# them to the closest width this engine's column-width grid can hit.
$ws.Columns.Item(2).ColumnWidth = 17.666666666666668
$ws.Columns.Item(3).ColumnWidth = 15.833333333333334
$ws.Columns.Item(5).ColumnWidth = 12.666666666666666
$ws.Columns.Item(6).ColumnWidth = 9.333333333333334
$ws.Columns.Item(7).ColumnWidth = 13.5
$ws.Columns.Item(8).ColumnWidth = 18.5
$ws.Columns.Item(9).ColumnWidth = 17.333333333333332
$ws.Columns.Item(10).ColumnWidth = 14.166666666666666
